$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 149, shifting rows 149:208 down to 150:209
$ws.Rows("149:149").Insert()

# Populate the newly inserted row 149 with the new weekly observation
$ws.Range("A149").Value = 7
$ws.Range("B149").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C149").Value = "Ñuble"
$ws.Range("D149").Value = 44704
$ws.Range("E149").Value = 16
$ws.Range("F149").Value = 100112017
$ws.Range("G149").Value = "Apio"
$ws.Range("H149").Value = "Americana (o)"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 100
$ws.Range("K149").Value = 7500
$ws.Range("L149").Value = 8000
$ws.Range("M149").Value = 7750
$ws.Range("N149").Value = "$/docena de matas"
$ws.Range("O149").Value = "Provincia del Elquí"
$ws.Range("P149").Value = 1292
$ws.Range("Q149").Value = 6
$ws.Range("R149").Value = "Hortaliza"

# Keep the date column's datetime display format consistent with the rest of column D
$ws.Range("D149").NumberFormat = $ws.Range("D150").NumberFormat
